$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with corrected analytics values
$ws.Range("A2").Value = 45880
$ws.Range("C2").Value = 5
$ws.Range("E2").Value = 4
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 2.75
$ws.Range("I2").Value = 0.39
$ws.Range("J2").Value = 71.40000000000001
$ws.Range("K2").Value = 28.59999999999999

# Remove row 3 entirely (no longer part of the dataset)
$ws.Rows("3").Delete()
